# Generate Report for Archive
#
# 1. The "Status" column value "Ready for handoff" becomes "In Translation"
#    on every sheet that shows it (Overview's per-language status columns,
#    plus each language sheet's own Status column).
# 2. The now-shorter status text causes the Status column(s) to re-autosize
#    narrower (width 17.2159881591797 -> 13.4101845877511) on every sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn status) and F (de-de status), rows 2-3
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: column C (Status), rows 2-3
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# de-de sheet: column C (Status), rows 2-3
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Resize the now-narrower Status columns to match the updated content width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
